# Update the "Product Search" sheet's PriceInResults / PriceInDetails /
# PriceInShoppingCart sample row (D2:F2) from "119.99" to "51.99".
#
# The target value must remain text (it is stored as a shared string in the
# workbook, like the other sample-data cells in that row) rather than being
# coerced into a Number cell. Assigning the numeric-looking string directly
# via .Value/.Value2/.Formula gets auto-coerced to a Number by Excel, and
# prefixing with a leading apostrophe forces Text but also stamps the cell
# with a quote-prefix flag (changing its style). Routing the text through a
# scratch cell's text-producing formula and a values-only paste keeps the
# cell's existing style/format untouched while still landing a Text value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Product Search")
$ws.Activate()

$scratch = $ws.Range("Z1")
$scratch.Formula = '="51.99"'
$scratch.Copy()

$ws.Range("D2").PasteSpecial(-4163)
$ws.Range("E2").PasteSpecial(-4163)
$ws.Range("F2").PasteSpecial(-4163)

$scratch.Clear()
$excel.CutCopyMode = $false
